# Correct column names used in workbook 2020
#
# The lookup-table sheets (Intercepts, Age, Gender, AdmitScore, AdmitPain,
# Payer, TreatmentType, Duration) all share one identical header row of
# column ids. Those ids were renamed from the old "IROMS" survey labels to
# the new "<MEASURE>.prob.<kind>" naming convention. Because every sheet
# points at the very same shared-string entries, updating the text once per
# distinct value (reusing the same literal string) lets the engine dedupe
# back onto the original shared-string slots, and the VLOOKUP-driven cache
# on "Example Calculation" picks up the rename automatically on recalc.

$wb = $excel.ActiveWorkbook

# old header (by column letter on the lookup sheets) -> new header text
$headerMap = [ordered]@{
    "B1" = "ODI.prob.mcd"
    "C1" = "ODI.prob.pain"
    "D1" = "NECK.prob.mcd"
    "E1" = "NECK.prob.pain"
    "F1" = "KNEE.prob.mcd"
    "G1" = "KNEE.prob.pain"
    "H1" = "LEFS.prob.mcd"
    "I1" = "LEFS.prob.pain"
    "J1" = "DASH.prob.mcd"
    "K1" = "DASH.prob.pain"
}

$lookupSheets = @("Intercepts", "Age", "Gender", "AdmitScore", "AdmitPain", "Payer", "TreatmentType")
foreach ($sheetName in $lookupSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $headerMap.Keys) {
        $ws.Range($addr).Value = $headerMap[$addr]
    }
}

# "Duration" carries the same header values, shifted two columns to the
# right (columns D:M instead of B:K) because A:C hold the duration bucket.
$durationHeaderMap = [ordered]@{
    "D1" = "ODI.prob.mcd"
    "E1" = "ODI.prob.pain"
    "F1" = "NECK.prob.mcd"
    "G1" = "NECK.prob.pain"
    "H1" = "KNEE.prob.mcd"
    "I1" = "KNEE.prob.pain"
    "J1" = "LEFS.prob.mcd"
    "K1" = "LEFS.prob.pain"
    "L1" = "DASH.prob.mcd"
    "M1" = "DASH.prob.pain"
}

$wsDuration = $wb.Worksheets.Item("Duration")
foreach ($addr in $durationHeaderMap.Keys) {
    $wsDuration.Range($addr).Value = $durationHeaderMap[$addr]
}

# Replay the navigation trail left behind in the saved view state: the user
# browsed through Age and Duration (leaving a new selection on each) before
# finally landing back on Intercepts, which ends up as the active tab.
$wsAge = $wb.Worksheets.Item("Age")
$wsAge.Activate()
$wsAge.Range("E19").Select()

$wsDuration.Activate()
$wsDuration.Range("D1:M1").Select()

$wsIntercepts = $wb.Worksheets.Item("Intercepts")
$wsIntercepts.Activate()
$wsIntercepts.Range("L37").Select()
